$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff removes two data rows from the "Export" sheet:
#   - account 008028807 / RAFAEL / 62000
#   - account 004948033 / GUILHERME / 11466.9
# Locate each row by its account number (column A) and delete the
# entire row, letting the remaining rows shift up naturally.

$row1 = $ws.Cells.Find("008028807")
$row1.EntireRow.Delete()

$row2 = $ws.Cells.Find("004948033")
$row2.EntireRow.Delete()
